$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.943671703338623
$ws.Range("B1").Value = 1.317521810531616
$ws.Range("C1").Value = 1.899257659912109
$ws.Range("D1").Value = 5.277941226959229
$ws.Range("E1").Value = 1.925993323326111
